$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 448
$ws1.Range("F5").Value = 30
$ws1.Range("F6").Value = 17
$ws1.Range("F7").Value = 253
$ws1.Range("F8").Value = 14095
$ws1.Range("F9").Value = 121
$ws1.Range("F10").Value = 101
$ws1.Range("F11").Value = 5669
$ws1.Range("F12").Value = 580
$ws1.Range("F16").Value = 1227
$ws1.Range("F17").Value = 2
$ws1.Range("F20").Value = 768
$ws1.Range("F21").Value = 2913
$ws1.Range("F23").Value = 10462
$ws1.Range("F25").Value = 42
$ws1.Range("F26").Value = 61

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 12

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 448
$ws4.Range("F6").Value = 30
$ws4.Range("F7").Value = 17
$ws4.Range("F8").Value = 253
$ws4.Range("F9").Value = 14095
$ws4.Range("F10").Value = 121
$ws4.Range("F11").Value = 101
$ws4.Range("F12").Value = 5669
$ws4.Range("F13").Value = 580
$ws4.Range("F17").Value = 1227
$ws4.Range("F18").Value = 2
$ws4.Range("F21").Value = 768
$ws4.Range("F22").Value = 2913
$ws4.Range("F24").Value = 12
$ws4.Range("F25").Value = 10462
$ws4.Range("F27").Value = 42
$ws4.Range("F28").Value = 61
